{"js": "// Update the worksheet date heading and regenerate every arithmetic\n// problem in the practice table (20 rows x 5 columns) with the new\n// values from the \"c986bee\" output refresh.\n//\n// The table's row/column layout (20x5 = 100 cells) is unchanged by\n// this edit -- only the text inside each cell (and the date line)\n// is replaced -- so we can do this as a single values-grid + one\n// paragraph text replace, which also preserves per-run formatting\n// (rFonts / sz) already present in each cell and in the heading run.\n\nconst newDate = \"2024-10-26 Saturday\";\n\nconst newValues = [\n  [\"75-17=\", \"19+25=\", \"70-22=\", \"8+78=\", \"29+56=\"],\n  [\"85-26=\", \"43+49=\", \"27+24=\", \"16+56=\", \"51-49=\"],\n  [\"71-33=\", \"19+62=\", \"62-36=\", \"16+68=\", \"46+37=\"],\n  [\"45+16=\", \"75-59=\", \"6+65=\", \"87-79=\", \"38+37=\"],\n  [\"39+18=\", \"12-7=\", \"64+7=\", \"52-33=\", \"82-63=\"],\n  [\"2+79=\", \"18+18=\", \"37-8=\", \"67+5=\", \"12+29=\"],\n  [\"9+12=\", \"43+49=\", \"88-59=\", \"17+34=\", \"21-8=\"],\n  [\"84-27=\", \"47-38=\", \"91-13=\", \"20-14=\", \"27+67=\"],\n  [\"64-18=\", \"93-75=\", \"43-36=\", \"48-19=\", \"23-16=\"],\n  [\"71-34=\", \"39+47=\", \"10-9=\", \"78-49=\", \"95-6=\"],\n  [\"22+49=\", \"88+9=\", \"43+9=\", \"53+18=\", \"12+19=\"],\n  [\"53-25=\", \"62-7=\", \"48+17=\", \"38+49=\", \"96-27=\"],\n  [\"84-16=\", \"86-49=\", \"43-7=\", \"83-25=\", \"60-26=\"],\n  [\"48+4=\", \"63-18=\", \"75+16=\", \"46+48=\", \"7+66=\"],\n  [\"12-4=\", \"16+79=\", \"17+46=\", \"58+9=\", \"36+6=\"],\n  [\"46+19=\", \"58-19=\", \"61-44=\", \"28+37=\", \"81-48=\"],\n  [\"71-68=\", \"15-6=\", \"19+67=\", \"27+45=\", \"5+77=\"],\n  [\"12+59=\", \"39+5=\", \"29+15=\", \"50-16=\", \"8+76=\"],\n  [\"49+2=\", \"39+26=\", \"55-9=\", \"18+78=\", \"23+69=\"],\n  [\"60-51=\", \"81-36=\", \"87-29=\", \"37+58=\", \"25+36=\"],\n];\n\nconst body = context.document.body;\n\n// 1) Update the date/day heading paragraph (first paragraph of the body),\n//    replacing its text in place so the existing run formatting\n//    (Arial, sz 30) is kept.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst headingParagraph = paragraphs.items[0];\nheadingParagraph.insertText(newDate, Word.InsertLocation.replace);\n\n// 2) Update every cell of the practice-problems table in one shot.\n//    Assigning `.values` rewrites each cell's text while leaving the\n//    table/row/cell structure (and each cell's run formatting) intact.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and regenerate every arithmetic\n# problem in the practice table (20 rows x 5 columns) with the new\n# values from the \"c986bee\" output refresh.\n#\n# The table's row/column layout (20x5 = 100 cells) is unchanged by\n# this edit -- only the text inside each cell (and the date line)\n# is replaced -- so we just overwrite each Range.Text in place, which\n# preserves the existing per-run formatting (rFonts / sz) already on\n# each cell/paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/day heading (first paragraph of the document).\n$d.Paragraphs.Item(1).Range.Text = \"2024-10-26 Saturday\"\n\n# 2) New values for every row of the practice-problems table, in\n#    reading order (row 1..20, column 1..5).\n$newValues = @(\n    @(\"75-17=\", \"19+25=\", \"70-22=\", \"8+78=\", \"29+56=\"),\n    @(\"85-26=\", \"43+49=\", \"27+24=\", \"16+56=\", \"51-49=\"),\n    @(\"71-33=\", \"19+62=\", \"62-36=\", \"16+68=\", \"46+37=\"),\n    @(\"45+16=\", \"75-59=\", \"6+65=\", \"87-79=\", \"38+37=\"),\n    @(\"39+18=\", \"12-7=\", \"64+7=\", \"52-33=\", \"82-63=\"),\n    @(\"2+79=\", \"18+18=\", \"37-8=\", \"67+5=\", \"12+29=\"),\n    @(\"9+12=\", \"43+49=\", \"88-59=\", \"17+34=\", \"21-8=\"),\n    @(\"84-27=\", \"47-38=\", \"91-13=\", \"20-14=\", \"27+67=\"),\n    @(\"64-18=\", \"93-75=\", \"43-36=\", \"48-19=\", \"23-16=\"),\n    @(\"71-34=\", \"39+47=\", \"10-9=\", \"78-49=\", \"95-6=\"),\n    @(\"22+49=\", \"88+9=\", \"43+9=\", \"53+18=\", \"12+19=\"),\n    @(\"53-25=\", \"62-7=\", \"48+17=\", \"38+49=\", \"96-27=\"),\n    @(\"84-16=\", \"86-49=\", \"43-7=\", \"83-25=\", \"60-26=\"),\n    @(\"48+4=\", \"63-18=\", \"75+16=\", \"46+48=\", \"7+66=\"),\n    @(\"12-4=\", \"16+79=\", \"17+46=\", \"58+9=\", \"36+6=\"),\n    @(\"46+19=\", \"58-19=\", \"61-44=\", \"28+37=\", \"81-48=\"),\n    @(\"71-68=\", \"15-6=\", \"19+67=\", \"27+45=\", \"5+77=\"),\n    @(\"12+59=\", \"39+5=\", \"29+15=\", \"50-16=\", \"8+76=\"),\n    @(\"49+2=\", \"39+26=\", \"55-9=\", \"18+78=\", \"23+69=\"),\n    @(\"60-51=\", \"81-36=\", \"87-29=\", \"37+58=\", \"25+36=\")\n)\n\n$tbl = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
